# Update the "想去人数" (number of people wanting to go) column F values
# for the 4 worksheets in the workbook, matching the gh-pages data refresh.

$wb = $excel.ActiveWorkbook

# Map: sheet name -> @{ row = newValue }
$updates = @{
    "展览"     = @{ 3 = 1747; 5 = 454; 8 = 1191; 9 = 340; 12 = 683; 13 = 184; 18 = 2922; 25 = 19; 26 = 5274; 31 = 309; 32 = 1096; 34 = 55 }
    "演出"     = @{ 4 = 1120; 14 = 609; 24 = 315; 26 = 3929; 30 = 198; 33 = 162 }
    "本地生活" = @{ 5 = 2454; 6 = 1041 }
    "全部类型" = @{ 4 = 2454; 5 = 1747; 6 = 1041; 11 = 454; 14 = 1191; 15 = 340; 17 = 683; 18 = 1120; 19 = 1120; 23 = 2922; 30 = 5274; 33 = 609; 34 = 609; 38 = 309; 44 = 315; 45 = 315; 47 = 1096; 48 = 198; 49 = 162; 50 = 55 }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Range("F$row").Value = $rows[$row]
    }
}
